# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to Leve profit columns (H-N) across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, as scraped from
# the scheduled "Omega_Profits" data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6000
$ws.Range("I40").Value = 5250
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 5250
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -5075
$ws.Range("N40").Value = -7850
$ws.Range("H51").Value = 11130.305
$ws.Range("I51").Value = 1499
$ws.Range("J51").Value = 12047.571
$ws.Range("K51").Value = 1499
$ws.Range("L51").Value = 12047.571
$ws.Range("M51").Value = -1015
$ws.Range("N51").Value = -13015.571
$ws.Range("H74").Value = 9068.263000000001
$ws.Range("I74").Value = 9076.294
$ws.Range("K74").Value = 9076.294
$ws.Range("M74").Value = -8140.294
$ws.Range("H77").Value = 9068.263000000001
$ws.Range("I77").Value = 9076.294
$ws.Range("K77").Value = 45381.47
$ws.Range("M77").Value = -40701.47
$ws.Range("H87").Value = 218536.75
$ws.Range("J87").Value = 218536.75
$ws.Range("L87").Value = 218536.75
$ws.Range("N87").Value = -221032.75
$ws.Range("H90").Value = 218536.75
$ws.Range("J90").Value = 218536.75
$ws.Range("L90").Value = 655610.25
$ws.Range("N90").Value = -668090.25
$ws.Range("H111").Value = 24437.25
$ws.Range("J111").Value = 24437.25
$ws.Range("L111").Value = 73311.75
$ws.Range("N111").Value = -79445.75
$ws.Range("H112").Value = 6368.8
$ws.Range("J112").Value = 6625.0527
$ws.Range("L112").Value = 19875.1581
$ws.Range("N112").Value = -22091.1581
$ws.Range("H137").Value = 2070.5386
$ws.Range("I137").Value = 1780.3
$ws.Range("J137").Value = 3038
$ws.Range("K137").Value = 5340.9
$ws.Range("L137").Value = 9114
$ws.Range("M137").Value = -2790.9
$ws.Range("N137").Value = -14214

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 4748
$ws.Range("I16").Value = 4748
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4748
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -4461
$ws.Range("H74").Value = 2209.4285
$ws.Range("I74").Value = 2189.422
$ws.Range("K74").Value = 2189.422
$ws.Range("M74").Value = -1315.422
$ws.Range("H77").Value = 2209.4285
$ws.Range("I77").Value = 2189.422
$ws.Range("K77").Value = 10947.11
$ws.Range("M77").Value = -6579.110000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2673.92
$ws.Range("I107").Value = 1991.7646
$ws.Range("K107").Value = 1991.7646
$ws.Range("M107").Value = -71.76459999999997
$ws.Range("H134").Value = 3107.8223
$ws.Range("I134").Value = 2673.525
$ws.Range("J134").Value = 6582.2
$ws.Range("K134").Value = 8020.575000000001
$ws.Range("L134").Value = 19746.6
$ws.Range("M134").Value = -5485.575000000001
$ws.Range("N134").Value = -24816.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 19519.312
$ws.Range("I16").Value = 12929.615
$ws.Range("J16").Value = 48074.668
$ws.Range("K16").Value = 12929.615
$ws.Range("L16").Value = 48074.668
$ws.Range("M16").Value = -12642.615
$ws.Range("N16").Value = -48648.668
$ws.Range("H31").Value = 6546.1904
$ws.Range("I31").Value = 7702.9614
$ws.Range("J31").Value = 4666.4375
$ws.Range("K31").Value = 7702.9614
$ws.Range("L31").Value = 4666.4375
$ws.Range("M31").Value = -7407.9614
$ws.Range("N31").Value = -5256.4375
$ws.Range("H34").Value = 6546.1904
$ws.Range("I34").Value = 7702.9614
$ws.Range("J34").Value = 4666.4375
$ws.Range("K34").Value = 7702.9614
$ws.Range("L34").Value = 4666.4375
$ws.Range("M34").Value = -7500.9614
$ws.Range("N34").Value = -5070.4375
$ws.Range("H36").Value = 5000
$ws.Range("I36").Value = 5000
$ws.Range("K36").Value = 5000
$ws.Range("M36").Value = -4612
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4840
$ws.Range("H58").Value = 3241.3076
$ws.Range("I58").Value = 3411.4167
$ws.Range("K58").Value = 3411.4167
$ws.Range("M58").Value = -3208.4167
$ws.Range("H113").Value = 19519.312
$ws.Range("I113").Value = 12929.615
$ws.Range("J113").Value = 48074.668
$ws.Range("K113").Value = 12929.615
$ws.Range("L113").Value = 48074.668
$ws.Range("M113").Value = -10759.615
$ws.Range("N113").Value = -52414.668
$ws.Range("H134").Value = 728.3889
$ws.Range("I134").Value = 700.9677
$ws.Range("K134").Value = 2102.9031
$ws.Range("M134").Value = 432.0969
$ws.Range("H136").Value = 3241.3076
$ws.Range("I136").Value = 3411.4167
$ws.Range("K136").Value = 10234.2501
$ws.Range("M136").Value = -7684.250100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 10331.333
$ws.Range("I64").Value = 9997.5
$ws.Range("K64").Value = 29992.5
$ws.Range("M64").Value = -29722.5
$ws.Range("H67").Value = 10331.333
$ws.Range("I67").Value = 9997.5
$ws.Range("K67").Value = 29992.5
$ws.Range("M67").Value = -29056.5
$ws.Range("H68").Value = 1429.4
$ws.Range("I68").Value = 749
$ws.Range("K68").Value = 2247
$ws.Range("M68").Value = -1436
$ws.Range("H71").Value = 1429.4
$ws.Range("I71").Value = 749
$ws.Range("K71").Value = 6741
$ws.Range("M71").Value = -2685
$ws.Range("H80").Value = 500
$ws.Range("I80").Value = 500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -564
$ws.Range("H81").Value = 5001
$ws.Range("J81").Value = 5001
$ws.Range("L81").Value = 15003
$ws.Range("N81").Value = -17249
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H83").Value = 500
$ws.Range("I83").Value = 500
$ws.Range("K83").Value = 4500
$ws.Range("M83").Value = 180
$ws.Range("H84").Value = 5001
$ws.Range("J84").Value = 5001
$ws.Range("L84").Value = 45009
$ws.Range("N84").Value = -56241
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H131").Value = 2043.9412
$ws.Range("I131").Value = 1293.625
$ws.Range("J131").Value = 2710.889
$ws.Range("K131").Value = 3880.875
$ws.Range("L131").Value = 8132.667
$ws.Range("M131").Value = 1159.125
$ws.Range("N131").Value = -18212.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3997.25
$ws.Range("I122").Value = 3996.6667
$ws.Range("K122").Value = 11990.0001
$ws.Range("M122").Value = -9540.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 28833.592
$ws.Range("I7").Value = 24717.87
$ws.Range("K7").Value = 24717.87
$ws.Range("M7").Value = -24605.87
$ws.Range("H17").Value = 11333.333
$ws.Range("I17").Value = 11333.333
$ws.Range("K17").Value = 11333.333
$ws.Range("M17").Value = -11163.333
$ws.Range("H56").Value = 15000
$ws.Range("I56").Value = 15000
$ws.Range("K56").Value = 15000
$ws.Range("M56").Value = -14309
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H109").Value = 358966.66
$ws.Range("I109").Value = 10000
$ws.Range("J109").Value = 533450
$ws.Range("K109").Value = 10000
$ws.Range("L109").Value = 533450
$ws.Range("M109").Value = -8613
$ws.Range("N109").Value = -536224
$ws.Range("H126").Value = 28833.592
$ws.Range("I126").Value = 24717.87
$ws.Range("K126").Value = 74153.61
$ws.Range("M126").Value = -71683.61
$ws.Range("H136").Value = 2871
$ws.Range("I136").Value = 2866.1667
$ws.Range("K136").Value = 8598.500100000001
$ws.Range("M136").Value = -6048.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = 0
$ws.Range("H81").Value = 2230.3845
$ws.Range("I81").Value = 2100
$ws.Range("J81").Value = 2665
$ws.Range("K81").Value = 4200
$ws.Range("L81").Value = 5330
$ws.Range("M81").Value = -3139
$ws.Range("N81").Value = -7452
$ws.Range("H84").Value = 2230.3845
$ws.Range("I84").Value = 2100
$ws.Range("J84").Value = 2665
$ws.Range("K84").Value = 21000
$ws.Range("L84").Value = 26650
$ws.Range("M84").Value = -15696
$ws.Range("N84").Value = -37258
$ws.Range("H96").Value = 125002670
$ws.Range("I96").Value = 166668580
$ws.Range("K96").Value = 166668580
$ws.Range("M96").Value = -166667207
$ws.Range("H113").Value = 1028.1666
$ws.Range("I113").Value = 947.75
$ws.Range("J113").Value = 1189
$ws.Range("K113").Value = 2843.25
$ws.Range("L113").Value = 3567
$ws.Range("M113").Value = -673.25
$ws.Range("N113").Value = -7907
$ws.Range("H122").Value = 3719.9285
$ws.Range("I122").Value = 3239.6667
$ws.Range("J122").Value = 4080.125
$ws.Range("K122").Value = 9719.000100000001
$ws.Range("L122").Value = 12240.375
$ws.Range("M122").Value = -7269.000100000001
$ws.Range("N122").Value = -17140.375
$ws.Range("H132").Value = 4238.75
$ws.Range("I132").Value = 4238.75
$ws.Range("K132").Value = 12716.25
$ws.Range("M132").Value = -10186.25
$ws.Range("H136").Value = 3852.8572
$ws.Range("I136").Value = 3078.3333
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 9234.999899999999
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -6684.999899999999
$ws.Range("N136").Value = -30600

Write-Host "Applied Omega_Profits updates across 8 sheets."
